$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.034.54"
$ws.Range("E2").Value = "  -2.05%  "
$ws.Range("D3").Value = "3.124.36"
$ws.Range("E3").Value = "  -0.47%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "'593.63"
$ws.Range("E5").Value = "  -2.55%  "
$ws.Range("D6").Value = "'136.37"
$ws.Range("E6").Value = "  -5.14%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "3.118.39"
$ws.Range("E8").Value = "  -0.55%  "
$ws.Range("E9").Value = "  -2.72%  "
$ws.Range("D10").Value = "'0.146"
$ws.Range("E10").Value = "  -3.19%  "
$ws.Range("D11").Value = "'5.20"
$ws.Range("E11").Value = "  -2.95%  "
$ws.Range("D12").Value = "'0.457"
$ws.Range("E12").Value = "  -3.46%  "
$ws.Range("D13").Value = "'0.0000247"
$ws.Range("E13").Value = "  -3.13%  "
$ws.Range("D14").Value = "'34.19"
$ws.Range("E14").Value = "  -3.43%  "
$ws.Range("D15").Value = "3.638.35"
$ws.Range("E15").Value = "  -0.58%  "
$ws.Range("E16").Value = "  +2.24%  "
$ws.Range("D17").Value = "3.135.17"
$ws.Range("E17").Value = "  -0.32%  "
$ws.Range("D18").Value = "63.003.09"
$ws.Range("E18").Value = "  -2.09%  "
$ws.Range("D19").Value = "'6.70"
$ws.Range("E19").Value = "  -2.51%  "
$ws.Range("D20").Value = "'474.45"
$ws.Range("E20").Value = "  -0.50%  "
$ws.Range("D21").Value = "'14.24"
$ws.Range("E21").Value = "  -3.60%  "
$ws.Range("D22").Value = "'0.699"
$ws.Range("E22").Value = "  -2.78%  "
$ws.Range("D23").Value = "'7.72"
$ws.Range("E23").Value = "  -1.18%  "
$ws.Range("D24").Value = "'86.98"
$ws.Range("E24").Value = "  +1.54%  "
$ws.Range("D25").Value = "'13.03"
$ws.Range("E25").Value = "  -4.13%  "
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("E27").Value = "  -2.07%  "
$ws.Range("D28").Value = "'7.17"
$ws.Range("E28").Value = "  -2.83%  "
$ws.Range("D29").Value = "'7.93"
$ws.Range("E29").Value = "  -6.25%  "
$ws.Range("D30").Value = "'2.05"
$ws.Range("E30").Value = "  -0.41%  "
$ws.Range("E31").Value = "  -0.11%  "
$ws.Range("D32").Value = "'26.70"
$ws.Range("E32").Value = "  +0.09%  "
$ws.Range("D33").Value = "'0.107"
$ws.Range("E33").Value = "  -8.53%  "
$ws.Range("D34").Value = "'2.53"
$ws.Range("E34").Value = "  -3.77%  "
$ws.Range("E35").Value = "  -2.69%  "
$ws.Range("D36").Value = "'5.82"
$ws.Range("E36").Value = "  -2.22%  "
$ws.Range("D37").Value = "'52.07"
$ws.Range("E37").Value = "  -1.24%  "
$ws.Range("D38").Value = "0.0₃0703"
$ws.Range("E38").Value = "  -5.40%  "
$ws.Range("D39").Value = "'0.0388"
$ws.Range("E39").Value = "  -1.64%  "
$ws.Range("D40").Value = "'420.68"
$ws.Range("E40").Value = "  -6.50%  "
$ws.Range("D41").Value = "'8.25"
$ws.Range("E41").Value = "  -0.87%  "
$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").Value = "'2.69"
$ws.Range("E42").Value = "  -10.03%  "
$ws.Range("D43").Value = "2.883.19"
$ws.Range("E43").Value = "  +0.09%  "
$ws.Range("B44").Value = "Kaspa"
$ws.Range("C44").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D44").Value = "'0.113"
$ws.Range("E44").Value = "  -4.76%  "
$ws.Range("D45").Value = "'0.262"
$ws.Range("E45").Value = "  -0.02%  "
$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").Value = "'2.13"
$ws.Range("E46").Value = "  -4.61%  "
$ws.Range("B47").Value = "USDe"
$ws.Range("C47").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D47").Value = "'0.999"
$ws.Range("E47").Value = "  -0.06%  "
$ws.Range("D48").Value = "'25.74"
$ws.Range("E48").Value = "  -2.66%  "
$ws.Range("E49").Value = "  -1.02%  "
$ws.Range("D50").Value = "'2.28"
$ws.Range("E50").Value = "  -5.88%  "
$ws.Range("D51").Value = "'119.61"
$ws.Range("E51").Value = "  -1.45%  "
